$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster FAPs -> ECs (Target cluster D2 stays ECs); refreshed TPM-derived metrics
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nmu"
$ws.Range("C2").Value = "Nmur1"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.147769
$ws.Range("H2").Value = 0.443307
$ws.Range("I2").Value = 0.5363704999673319
$ws.Range("J2").Value = 0.5363704999673319
$ws.Range("M2").Value = 0.02360366666666667
$ws.Range("N2").Value = 0.070811
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.003487890219666667
$ws.Range("R2").Value = 0.031391011977
$ws.Range("S2").Value = 0.5363704999673319
$ws.Range("T2").Value = 0.5363704999673319

# Row 3: Sending cluster stays FAPs; Target cluster Inflammatory-Mac -> ECs; refreshed TPM-derived metrics
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Nmu"
$ws.Range("C3").Value = "Nmur1"
$ws.Range("D3").Value = "ECs"
$ws.Range("I3").Value = 0.4636295000326681
$ws.Range("J3").Value = 0.4636295000326681
$ws.Range("M3").Value = 0.02360366666666667
$ws.Range("N3").Value = 0.070811
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.003014872739666666
$ws.Range("R3").Value = 0.027133854657
$ws.Range("S3").Value = 0.4636295000326681
$ws.Range("T3").Value = 0.4636295000326681
